$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1. Update / fill the data for the new "Akka Distributed (2 physical
#    machines)" column (column E) on Sheet1.
# ------------------------------------------------------------------
$ws.Range("E3").Value = "Akka Distributed" + [char]10 + "(2 pysical machines)"
$ws.Range("E5").Value = 185226
$ws.Range("E6").Value = 370398
$ws.Range("E7").Value = 556172
$ws.Range("E8").Value = 741398

# ------------------------------------------------------------------
# 2. Update the selection / active cell shown on Sheet1 (cosmetic, but
#    matches the author's recorded selection at save-time).
# ------------------------------------------------------------------
$ws.Range("A1").Select() | Out-Null
$ws.Range("A2:E2").Select() | Out-Null

# ------------------------------------------------------------------
# 3. Chart 1 - "Processing Time Comparison (1 million passwords)"
#    (bar3D chart). Extend the single series from B:D to B:E so the
#    new Akka (2 machines) bar shows up, and drop the hard-coded
#    series title text (it becomes un-named, same as the source file).
# ------------------------------------------------------------------
$chart1 = $ws.ChartObjects().Item(1).Chart
$ser1 = $chart1.SeriesCollection().Item(1)
$ser1.Delete() | Out-Null
$newSer1 = $chart1.SeriesCollection().NewSeries()
$newSer1.Formula = "=SERIES(,Sheet1!`$B`$3:`$E`$3,Sheet1!`$B`$5:`$E`$5,1)"
$newSer1.InvertIfNegative = $false

# ------------------------------------------------------------------
# 4. Chart 2 - the line chart. Add the 4th series for the new Akka (2
#    machines) column, following the same pattern as the 3 existing
#    series (name from row 3, categories from A4:A8, values from the
#    new column E4:E8).
# ------------------------------------------------------------------
$chart2 = $ws.ChartObjects().Item(2).Chart
$newSer2 = $chart2.SeriesCollection().NewSeries()
$newSer2.Formula = "=SERIES(Sheet1!`$E`$3,Sheet1!`$A`$4:`$A`$8,Sheet1!`$E`$4:`$E`$8,4)"
$newSer2.Smooth = $false
